$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G2").Value = "2016-08-19 19:12:25"
$dede.Range("H2").Value = "2016-08-19 19:12:25"
$zhcn.Range("H2").Value = "2016-08-19 19:12:15"
$zhcn.Range("K2").Value = "2016-08-19 19:12:47"
$dede.Range("K2").Value = "2016-08-19 19:12:54"
